$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.991.49"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "2.422.99"
$ws.Range("E3").Value = "  -2.68%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.54"
$ws.Range("E5").Value = "  -3.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.19"
$ws.Range("E6").Value = "  -2.93%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.496"
$ws.Range("E8").Value = "  -2.83%  "

$ws.Range("D9").Value = "2.423.58"
$ws.Range("E9").Value = "  -2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  -8.63%  "

$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("E12").Value = "  -6.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.73"
$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("D14").Value = "2.877.41"
$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("D15").Value = "68.069.93"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -5.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.00"
$ws.Range("E17").Value = "  -4.57%  "

$ws.Range("D18").Value = "2.461.49"
$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  -5.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.30"
$ws.Range("E20").Value = "  -2.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("E21").Value = "  -6.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.73"
$ws.Range("E22").Value = "  -3.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  -3.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.93"
$ws.Range("E25").Value = "  -4.85%  "

$ws.Range("D26").Value = "2.572.05"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.57"
$ws.Range("E27").Value = "  -8.00%  "

$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.93"
$ws.Range("E29").Value = "  -7.88%  "

$ws.Range("D30").Value = "0.0₃0793"
$ws.Range("E30").Value = "  -8.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  -7.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "428.62"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.10"
$ws.Range("E34").Value = "  -7.17%  "

$ws.Range("E35").Value = "  -6.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.66"
$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.99"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  -4.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.58"
$ws.Range("E40").Value = "  -2.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.297"
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.31"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.32"
$ws.Range("E43").Value = "  -5.21%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.06"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.42"
$ws.Range("E45").Value = "  -9.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("E46").Value = "  -8.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "129.35"
$ws.Range("E47").Value = "  -6.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.28"
$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.474"
$ws.Range("E50").Value = "  -5.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.550"
$ws.Range("E51").Value = "  -3.73%  "
